$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.282.07"
$ws.Range("E2").Value = "  +0.13%  "

$ws.Range("D3").Value = "3.569.64"
$ws.Range("E3").Value = "  +0.74%  "

$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.12%  "

$ws.Range("D7").Value = "3.568.10"
$ws.Range("E7").Value = "  +0.56%  "

$ws.Range("E8").Value = "  +0.22%  "

$ws.Range("E9").Value = "  +2.13%  "

$ws.Range("E10").Value = "  -0.15%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.83"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.77%  "

$ws.Range("E12").Value = "  -0.10%  "

$ws.Range("D13").Value = "4.173.22"
$ws.Range("E13").Value = "  +0.81%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000207"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.67%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "30.41"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.22%  "

$ws.Range("D16").Value = "3.544.43"
$ws.Range("E16").Value = "  +0.07%  "

$ws.Range("D17").Value = "66.302.88"
$ws.Range("E17").Value = "  -0.05%  "

$ws.Range("E18").Value = "  -0.55%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.33%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.24%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "430.96"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.15%  "

$ws.Range("E23").Value = "  +1.95%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.42%  "

$ws.Range("D25").Value = "3.710.74"
$ws.Range("E25").Value = "  +0.83%  "

$ws.Range("E26").Value = "  +0.08%  "

$ws.Range("E27").Value = "  -1.82%  "

$ws.Range("E28").Value = "  +1.01%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.14%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.93"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.97%  "

$ws.Range("E31").Value = "  -0.16%  "

$ws.Range("D32").Value = "3.562.54"
$ws.Range("E32").Value = "  +0.92%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "25.48"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.41%  "

$ws.Range("E34").Value = "  -2.42%  "

$ws.Range("E35").Value = "  -4.39%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.84"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.44%  "

$ws.Range("E38").Value = "  -1.66%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.62"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "175.52"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.58%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0852"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.27%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.22"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.50%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.889"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.69%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.93"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "46.04"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.62%  "

$ws.Range("E46").Value = "  -0.11%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.50"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.39%  "

$ws.Range("E48").Value = "  -2.24%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.12"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.76%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.76%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.30%  "
